$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin data (prices + 1h volume change) scraped on
# Sun Jul 23 19:38:17 UTC 2023. Some rows were re-ranked (12-14).

$ws.Range('D2').Value = '30.186.05'
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').Value = '1.893.96'
$ws.Range('E3').Value = '  +0.39%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7446'
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3176'
$ws.Range('E8').Value = '  +2.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07238'
$ws.Range('E9').Value = '  +1.84%  '
$ws.Range('E10').Value = '  -1.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08362'
$ws.Range('E11').Value = '  -1.47%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.954.51'
$ws.Range('E12').Value = '  +3.08%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7633'
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.448'
$ws.Range('E14').Value = '  +1.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.06'
$ws.Range('E15').Value = '  -0.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.189'
$ws.Range('E16').Value = '  +0.73%  '
$ws.Range('D17').Value = '30.245.74'
$ws.Range('E17').Value = '  +1.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '250.91'
$ws.Range('E18').Value = '  +3.21%  '
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007876'
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('D21').Value = '2.157.19'
$ws.Range('E21').Value = '  -0.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.008'
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.000'
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1588'
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.317'
$ws.Range('E26').Value = '  -0.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '164.30'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.81'
$ws.Range('E28').Value = '  +0.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.075'
$ws.Range('E29').Value = '  +2.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.480'
$ws.Range('E30').Value = '  -1.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.611'
$ws.Range('E31').Value = '  +3.26%  '
$ws.Range('E32').Value = '  +0.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.233'
$ws.Range('E33').Value = '  +3.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05408'
$ws.Range('E34').Value = '  +0.24%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.254'
$ws.Range('E35').Value = '  +1.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7680'
$ws.Range('E36').Value = '  +3.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9928'
$ws.Range('E37').Value = '  -0.93%  '
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('E39').Value = '  +2.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.771'
$ws.Range('E40').Value = '  +0.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4576'
$ws.Range('E41').Value = '  +2.86%  '
$ws.Range('D42').Value = '1.103.21'
$ws.Range('E42').Value = '  +1.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.095'
$ws.Range('E43').Value = '  +0.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.85'
$ws.Range('E44').Value = '  +0.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8710'
$ws.Range('E45').Value = '  +0.96%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '104.35'
$ws.Range('E46').Value = '  +1.73%  '
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('E48').Value = '  +0.77%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.636'
$ws.Range('E49').Value = '  -0.34%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.634'
$ws.Range('E50').Value = '  -0.50%  '
$ws.Range('D51').Value = '2.061.81'
$ws.Range('E51').Value = '  +0.30%  '
